# Test case 13 completed, is successfully running
# Add two new rows of effort-tracking data to the active sheet ("Effort R 1.0")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17: 20/06/2013, 0.5h regular effort, task description
$ws.Range("A17").Value = 41445
$ws.Range("B17").Value = 0.5
$ws.Range("D17").Value = "Implementation tc13_eventStates"

# Row 18: 21/06/2013, 2.25h additional effort, task description
$ws.Range("A18").Value = 41446
$ws.Range("C18").Value = 2.25
$ws.Range("D18").Value = "Successful completion of tc13"

# Update the visible selection to match the edited workbook state
[void]$ws.Range("E11").Select()
